$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate formatting (styles) from the last existing data row (109) to the new rows,
# so the new rows match the workbook's established look (bold/bordered index column,
# date-time format on the "data_partida" column, etc.)
$ws.Range("A109:V109").Copy() | Out-Null
$ws.Range("A110:V113").PasteSpecial(-4122) | Out-Null

# Column D ("temporada") holds a numeric-looking value ("2023") that must stay text,
# like the rest of the sheet - force a text number format before writing it so Excel
# does not auto-convert it to a real number.
$ws.Range("D110:D113").NumberFormat = "@"

# Row 110 (Indice 109)
$ws.Range("A110").Value = 109
$ws.Range("B110").Value = 'ecuador'
$ws.Range("C110").Value = 'liga-pro'
$ws.Range("D110").Value = '2023'
$ws.Range("E110").Value = 45257
$ws.Range("F110").Value = 'Gualaceo'
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 'Emelec'
$ws.Range("I110").Value = 2
$ws.Range("J110").Value = 3.41
$ws.Range("K110").Value = '20/11/2023 00:12'
$ws.Range("L110").Value = 2.62
$ws.Range("M110").Value = '26/11/2023 23:58'
$ws.Range("N110").Value = 3.32
$ws.Range("O110").Value = '20/11/2023 00:12'
$ws.Range("P110").Value = 3.33
$ws.Range("Q110").Value = '26/11/2023 23:58'
$ws.Range("R110").Value = 2.18
$ws.Range("S110").Value = '20/11/2023 00:12'
$ws.Range("T110").Value = 2.76
$ws.Range("U110").Value = '26/11/2023 23:58'
$ws.Range("V110").Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/gualaceo-emelec/8MLRzLQ2/'

# Row 111 (Indice 110)
$ws.Range("A111").Value = 110
$ws.Range("B111").Value = 'ecuador'
$ws.Range("C111").Value = 'liga-pro'
$ws.Range("D111").Value = '2023'
$ws.Range("E111").Value = 45257
$ws.Range("F111").Value = 'Barcelona SC'
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 'Guayaquil City'
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 1.28
$ws.Range("K111").Value = '20/11/2023 00:12'
$ws.Range("L111").Value = 1.45
$ws.Range("M111").Value = '26/11/2023 23:58'
$ws.Range("N111").Value = 5.6
$ws.Range("O111").Value = '20/11/2023 00:12'
$ws.Range("P111").Value = 4.42
$ws.Range("Q111").Value = '26/11/2023 23:58'
$ws.Range("R111").Value = 10.05
$ws.Range("S111").Value = '20/11/2023 00:12'
$ws.Range("T111").Value = 7.46
$ws.Range("U111").Value = '26/11/2023 23:58'
$ws.Range("V111").Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/barcelona-sc-guayaquil-city/nDKVZvt9/'

# Row 112 (Indice 111)
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = 'ecuador'
$ws.Range("C112").Value = 'liga-pro'
$ws.Range("D112").Value = '2023'
$ws.Range("E112").Value = 45257
$ws.Range("F112").Value = 'Cumbaya'
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 'LDU Quito'
$ws.Range("I112").Value = 2
$ws.Range("J112").Value = 5.71
$ws.Range("K112").Value = '20/11/2023 00:12'
$ws.Range("L112").Value = 9.789999999999999
$ws.Range("M112").Value = '26/11/2023 23:55'
$ws.Range("N112").Value = 4.01
$ws.Range("O112").Value = '20/11/2023 00:12'
$ws.Range("P112").Value = 4.53
$ws.Range("Q112").Value = '26/11/2023 23:55'
$ws.Range("R112").Value = 1.53
$ws.Range("S112").Value = '20/11/2023 00:12'
$ws.Range("T112").Value = 1.37
$ws.Range("U112").Value = '26/11/2023 23:55'
$ws.Range("V112").Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/cumbaya-ldu-quito/tzKZYbeF/'

# Row 113 (Indice 112)
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = 'ecuador'
$ws.Range("C113").Value = 'liga-pro'
$ws.Range("D113").Value = '2023'
$ws.Range("E113").Value = 45257
$ws.Range("F113").Value = 'Tecnico U.'
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 'Libertad'
$ws.Range("I113").Value = 1
$ws.Range("J113").Value = 1.63
$ws.Range("K113").Value = '20/11/2023 00:12'
$ws.Range("L113").Value = 1.55
$ws.Range("M113").Value = '26/11/2023 23:21'
$ws.Range("N113").Value = 3.75
$ws.Range("O113").Value = '20/11/2023 00:12'
$ws.Range("P113").Value = 4.04
$ws.Range("Q113").Value = '26/11/2023 23:21'
$ws.Range("R113").Value = 5.56
$ws.Range("S113").Value = '20/11/2023 00:12'
$ws.Range("T113").Value = 6.39
$ws.Range("U113").Value = '26/11/2023 23:21'
$ws.Range("V113").Value = 'https://www.betexplorer.com/football/ecuador/liga-pro/tecnico-u-libertad/ny3IQz3e/'
